$wb = $excel.ActiveWorkbook

$wsRecap = $wb.Worksheets.Item("Recap")
$wsColab = $wb.Worksheets.Item("Colab")
$wsListe = $wb.Worksheets.Item("Liste")

# --- Colab sheet: insert new collaborator row (Mehdi Tajmouati / MTI) ---
# Row 5 is currently an empty gap row between row 4 (AVE) and row 6 (blank
# styled footer row); inserting here pushes the footer row down to row 7.
$wsColab.Rows("5:5").Insert()
$wsColab.Range("B5").Value = "Mehdi Tajmouati"
$wsColab.Range("A5").Value = "MTI"

# --- Liste sheet: insert new client CASA, keeping the Client list sorted ---
# (only column G is part of the Client list; shift G5:G8 down into G6:G9
# and place CASA at G5, just before CGI)
$wsListe.Range("G9").Value = $wsListe.Range("G8").Value()
$wsListe.Range("G8").Value = $wsListe.Range("G7").Value()
$wsListe.Range("G7").Value = $wsListe.Range("G6").Value()
$wsListe.Range("G6").Value = $wsListe.Range("G5").Value()
$wsListe.Range("G5").Value = "CASA"

# --- Recap sheet: insert a fresh blank template row at row 8 (this shifts
# the old blank row 8 -> 9 and the footer border row 9 -> 10), then fill
# the now-vacated row 8 with the new Mehdi Tajmouati / CASA entry ---
$wsRecap.Rows("8:8").Insert()
$wsRecap.Range("A8").Value = "MTI"
$wsRecap.Range("B8").Value = "TALEND"
$wsRecap.Range("C8").Value = 3
$wsRecap.Range("D8").Value = "ETL"
$wsRecap.Range("E8").Value = 3
$wsRecap.Range("F8").Value = "CASA"
$wsRecap.Range("G8").Value = 18

# --- Defined names: extend Client and Colab ranges by one row each ---
$wb.Names.Item("Client").RefersTo = "=Liste!`$G`$2:`$G`$9"
$wb.Names.Item("Colab").RefersTo = "=Colab!`$A`$2:`$A`$6"

# --- Selections / active sheet ---
$wsRecap.Range("G9").Select()
$wsListe.Range("F10").Select()
$wsColab.Range("B5").Select()
$wsColab.Activate()
